$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing Code/Description/Definition columns (A/B/C) one column to the
# right (B/C/D) to make room for the new "Version" column in A. Using Copy (rather
# than re-assigning .Value) preserves the original shared-string cell typing and
# avoids touching cell styles.
$ws.Range("A1:C5").Copy($ws.Range("B1"))

# New header for the inserted column
$ws.Range("A1").Value = "Version"

# Fill the version value "1.0" for each data row. A plain .Value assignment of a
# numeric-looking string like "1.0" gets silently coerced into a number, so instead
# build it as a text-formula result in a helper range, then copy/paste-values it
# back in so it ends up stored as the shared string "1.0" rather than a number.
$helper = $ws.Range("F1:F4")
$helper.Formula = '="1.0"'
$helper.Copy()
$ws.Range("A2:A5").PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = 0
